$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-12 03:28:56"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-12 03:28:50"
$wsZhCn.Range("K2").Value = "2016-08-12 03:29:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-12 03:28:56"
$wsDeDe.Range("K2").Value = "2016-08-12 03:29:20"
